$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3776.625
$ws.Range("I64").Value = 3278.5264
$ws.Range("J64").Value = 4102.9653
$ws.Range("K64").Value = 3278.5264
$ws.Range("L64").Value = 4102.9653
$ws.Range("M64").Value = -3030.5264
$ws.Range("N64").Value = -4598.9653
$ws.Range("H67").Value = 3776.625
$ws.Range("I67").Value = 3278.5264
$ws.Range("J67").Value = 4102.9653
$ws.Range("K67").Value = 3278.5264
$ws.Range("L67").Value = 4102.9653
$ws.Range("M67").Value = -2420.5264
$ws.Range("N67").Value = -5818.9653
$ws.Range("H96").Value = 686.6923
$ws.Range("I96").Value = 316.85715
$ws.Range("J96").Value = 1118.1666
$ws.Range("K96").Value = 950.5714499999999
$ws.Range("L96").Value = 3354.4998
$ws.Range("M96").Value = 422.4285500000001
$ws.Range("N96").Value = -6100.4998
$ws.Range("H137").Value = 13758894
$ws.Range("I137").Value = 25452124
$ws.Range("J137").Value = 2152.1765
$ws.Range("K137").Value = 76356372
$ws.Range("L137").Value = 6456.529500000001
$ws.Range("M137").Value = -76353822
$ws.Range("N137").Value = -11556.5295
$ws.Range("H138").Value = 2473.23
$ws.Range("I138").Value = 874.75
$ws.Range("J138").Value = 2872.85
$ws.Range("K138").Value = 2624.25
$ws.Range("L138").Value = 8618.549999999999
$ws.Range("M138").Value = 2515.75
$ws.Range("N138").Value = -18898.55
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3026.1
$ws.Range("I61").Value = 2754.4614
$ws.Range("J61").Value = 3530.5715
$ws.Range("K61").Value = 2754.4614
$ws.Range("L61").Value = 3530.5715
$ws.Range("M61").Value = -2542.4614
$ws.Range("N61").Value = -3954.5715
$ws.Range("H74").Value = 2771.9167
$ws.Range("I74").Value = 1322.5
$ws.Range("J74").Value = 4221.3335
$ws.Range("K74").Value = 1322.5
$ws.Range("L74").Value = 4221.3335
$ws.Range("M74").Value = -448.5
$ws.Range("N74").Value = -5969.3335
$ws.Range("H77").Value = 2771.9167
$ws.Range("I77").Value = 1322.5
$ws.Range("J77").Value = 4221.3335
$ws.Range("K77").Value = 6612.5
$ws.Range("L77").Value = 21106.6675
$ws.Range("M77").Value = -2244.5
$ws.Range("N77").Value = -29842.6675
$ws.Range("H97").Value = 426.1875
$ws.Range("I97").Value = 388
$ws.Range("J97").Value = 999
$ws.Range("K97").Value = 388
$ws.Range("L97").Value = 999
$ws.Range("M97").Value = 108
$ws.Range("N97").Value = -1991
$ws.Range("H102").Value = 3065.125
$ws.Range("I102").Value = 3470
$ws.Range("J102").Value = 2822.2
$ws.Range("K102").Value = 3470
$ws.Range("L102").Value = 2822.2
$ws.Range("M102").Value = -1848
$ws.Range("N102").Value = -6066.2
$ws.Range("H132").Value = 1685.0294
$ws.Range("I132").Value = 932.13635
$ws.Range("J132").Value = 3065.3333
$ws.Range("K132").Value = 2796.40905
$ws.Range("L132").Value = 9195.999899999999
$ws.Range("M132").Value = -266.4090500000002
$ws.Range("N132").Value = -14255.9999
$ws.Range("H136").Value = 3026.1
$ws.Range("I136").Value = 2754.4614
$ws.Range("J136").Value = 3530.5715
$ws.Range("K136").Value = 8263.3842
$ws.Range("L136").Value = 10591.7145
$ws.Range("M136").Value = -5713.3842
$ws.Range("N136").Value = -15691.7145
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1745.5
$ws.Range("I86").Value = 1663.6666
$ws.Range("J86").Value = 1991
$ws.Range("K86").Value = 1663.6666
$ws.Range("L86").Value = 1991
$ws.Range("M86").Value = -540.6666
$ws.Range("N86").Value = -4237
$ws.Range("H89").Value = 1745.5
$ws.Range("I89").Value = 1663.6666
$ws.Range("J89").Value = 1991
$ws.Range("K89").Value = 8318.333000000001
$ws.Range("L89").Value = 9955
$ws.Range("M89").Value = -2702.333000000001
$ws.Range("N89").Value = -21187
$ws.Range("H94").Value = 278.57144
$ws.Range("I94").Value = 191.66667
$ws.Range("J94").Value = 800
$ws.Range("K94").Value = 191.66667
$ws.Range("L94").Value = 800
$ws.Range("M94").Value = 259.33333
$ws.Range("N94").Value = -1702
$ws.Range("H99").Value = 2013.1578
$ws.Range("I99").Value = 1523.4286
$ws.Range("J99").Value = 3384.4
$ws.Range("K99").Value = 1523.4286
$ws.Range("L99").Value = 3384.4
$ws.Range("M99").Value = -25.42859999999996
$ws.Range("N99").Value = -6380.4
$ws.Range("H105").Value = 1984.6471
$ws.Range("I105").Value = 1475.5834
$ws.Range("J105").Value = 3206.4
$ws.Range("K105").Value = 1475.5834
$ws.Range("L105").Value = 3206.4
$ws.Range("M105").Value = 271.4166
$ws.Range("N105").Value = -6700.4
$ws.Range("H134").Value = 2761.4194
$ws.Range("I134").Value = 1937.0952
$ws.Range("J134").Value = 4492.5
$ws.Range("K134").Value = 5811.2856
$ws.Range("L134").Value = 13477.5
$ws.Range("M134").Value = -3276.2856
$ws.Range("N134").Value = -18547.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2058.743
$ws.Range("I58").Value = 1411.174
$ws.Range("J58").Value = 3299.9167
$ws.Range("K58").Value = 1411.174
$ws.Range("L58").Value = 3299.9167
$ws.Range("M58").Value = -1208.174
$ws.Range("N58").Value = -3705.9167
$ws.Range("H134").Value = 3305.1667
$ws.Range("I134").Value = 1618.8
$ws.Range("J134").Value = 6115.778
$ws.Range("K134").Value = 4856.4
$ws.Range("L134").Value = 18347.334
$ws.Range("M134").Value = -2321.4
$ws.Range("N134").Value = -23417.334
$ws.Range("H136").Value = 2058.743
$ws.Range("I136").Value = 1411.174
$ws.Range("J136").Value = 3299.9167
$ws.Range("K136").Value = 4233.522
$ws.Range("L136").Value = 9899.750100000001
$ws.Range("M136").Value = -1683.522
$ws.Range("N136").Value = -14999.7501
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H93").Value = 5000
$ws.Range("J93").Value = 5000
$ws.Range("L93").Value = 15000
$ws.Range("N93").Value = -18744
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1049.0714
$ws.Range("I97").Value = 1049.0714
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1049.0714
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -553.0714
$ws.Range("N97").ClearContents()
$ws.Range("H132").Value = 7038.273
$ws.Range("I132").Value = 10753
$ws.Range("J132").Value = 4915.5713
$ws.Range("K132").Value = 32259
$ws.Range("L132").Value = 14746.7139
$ws.Range("M132").Value = -29729
$ws.Range("N132").Value = -19806.7139
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 14287481
$ws.Range("I100").Value = 1704.15
$ws.Range("J100").Value = 33335182
$ws.Range("K100").Value = 1704.15
$ws.Range("L100").Value = 33335182
$ws.Range("M100").Value = -1163.15
$ws.Range("N100").Value = -33336264
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 16702338
$ws.Range("I62").Value = 20041636
$ws.Range("J62").Value = 5850
$ws.Range("K62").Value = 20041636
$ws.Range("L62").Value = 5850
$ws.Range("M62").Value = -20041012
$ws.Range("N62").Value = -7098
$ws.Range("H65").Value = 16702338
$ws.Range("I65").Value = 20041636
$ws.Range("J65").Value = 5850
$ws.Range("K65").Value = 100208180
$ws.Range("L65").Value = 29250
$ws.Range("M65").Value = -100205060
$ws.Range("N65").Value = -35490
$ws.Range("H81").Value = 2491.6667
$ws.Range("I81").Value = 2088.889
$ws.Range("J81").Value = 3700
$ws.Range("K81").Value = 4177.778
$ws.Range("L81").Value = 7400
$ws.Range("M81").Value = -3116.778
$ws.Range("N81").Value = -9522
$ws.Range("H84").Value = 2491.6667
$ws.Range("I84").Value = 2088.889
$ws.Range("J84").Value = 3700
$ws.Range("K84").Value = 20888.89
$ws.Range("L84").Value = 37000
$ws.Range("M84").Value = -15584.89
$ws.Range("N84").Value = -47608
